$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# "name" column value for the pit-size row changed from "Perc2.97" to "pct_9m2"
$ws.Range("B8").Value = "pct_9m2"

# Update the active selection on the sheet to match the target state
$ws.Range("B9").Select()
